# Apply the "new .ttl from Google sheet has been generated" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) ConceptScheme / PREFIX base URI changed from the old m4m-dk-Test4
#    namespace to the new test3 namespace (two occurrences).
$ws.Range("B1").Value = "http://purl.org/test3/variables/"
$ws.Range("C3").Value = "http://purl.org/test3/variables/"

# 2) Fill in vocabulary metadata values that were previously left blank.
$ws.Range("B10").Value = "test"     # dct:title
$ws.Range("B11").Value = "test"     # dct:description
$ws.Range("B12").Value = "Hannah"   # dct:creator (first creator)

# 3) Insert a new row above the old row 13 (dct:rights) for a second
#    dct:creator entry, pushing every following row down by one.
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "dct:creator"
$ws.Range("B13").Value = "Minka"

# 4) Populate two of the previously-empty "vars:" placeholder rows (now at
#    rows 20 and 21 after the insert above) with real term data.
$ws.Range("A20").Value = "vars:test"
$ws.Range("B20").Value = "test"

$ws.Range("A21").Value = "vars:computerscientist"
$ws.Range("B21").Value = "computerscientist"
$ws.Range("E21").Value = "a person that knows stuff about computers"
